$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine how many data rows are present (header is row 1, data starts row 2).
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# For every session row whose Status (column I) is "Recorded", record who
# took attendance in the "Recorded By" column (column G).
for ($r = 2; $r -le $lastRow; $r++) {
    $status = $ws.Cells.Item($r, 9).Value2
    if ($status -eq "Recorded") {
        $ws.Cells.Item($r, 7).Value = "Miss Dina Nasr, Administrator"
    }
}

# Widen column G ("Recorded By") so the administrator's name fits/displays
# fully (stored column width of 31 characters).
$ws.Columns.Item(7).ColumnWidth = 31 - 5/6
